$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (D, E) for ownTeam / oppTeam, shifting old D:I -> F:K
$ws.Range("D1:E1").EntireColumn.Insert()

# Header row
$ws.Range("A1").Value = "venue"
$ws.Range("B1").Value = "date"
$ws.Range("C1").Value = "result"
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"
$ws.Range("F1").Value = "batsman"
$ws.Range("G1").Value = "totalRuns"
$ws.Range("H1").Value = "totalBalls"
$ws.Range("I1").Value = "total4s"
$ws.Range("J1").Value = "total6s"
$ws.Range("K1").Value = "sr"

# Clear old data rows before rewriting (also removes stale F:K values left after column insert)
$ws.Range("A2:K13").ClearContents()

# Data rows - write as Text (matches source formatting: numbers stored as text strings)
$data = @(
    @(" Abu Dhabi", " October 25 2020", "Royals won by 8 wickets (with 10 balls remaining)", "Mumbai Indians", "Rajasthan Royals", "Krunal Pandya ", "3", "4", "0", "0", "75.00"),
    @(" Abu Dhabi", " October 28 2020", "Mumbai won by 5 wickets (with 5 balls remaining)", "Mumbai Indians", "Royal Challengers Bangalore", "Krunal Pandya ", "10", "10", "1", "0", "100.00"),
    @(" Sharjah", " November 03 2020", "Sunrisers won by 10 wickets (with 17 balls remaining)", "Mumbai Indians", "Sunrisers Hyderabad", "Krunal Pandya ", "0", "3", "0", "0", "0.00"),
    @(" Dubai (DSC)", " November 05 2020", "Mumbai won by 57 runs", "Mumbai Indians", "Delhi Capitals", "Krunal Pandya ", "13", "10", "0", "1", "130.00"),
    @(" Dubai (DSC)", " November 10 2020", "Mumbai won by 5 wickets (with 8 balls remaining)", "Mumbai Indians", "Delhi Capitals", "Krunal Pandya ", "1", "1", "0", "0", "100.00"),
    @(" Abu Dhabi", " September 19 2020", "Super Kings won by 5 wickets (with 4 balls remaining)", "Mumbai Indians", "Chennai Super Kings", "Krunal Pandya ", "3", "3", "0", "0", "100.00"),
    @(" Abu Dhabi", " October 06 2020", "Mumbai won by 57 runs", "Mumbai Indians", "Rajasthan Royals", "Krunal Pandya ", "12", "17", "0", "1", "70.58"),
    @(" Dubai (DSC)", " October 18 2020", "Match tied (Kings XI won the one-over eliminator)", "Mumbai Indians", "Kings XI Punjab", "Krunal Pandya ", "34", "30", "4", "1", "113.33"),
    @(" Dubai (DSC)", " September 28 2020", "Match tied (RCB won the one-over eliminator)", "Mumbai Indians", "Royal Challengers Bangalore", "Krunal Pandya ", "0", "0", "0", "0", "-"),
    @(" Abu Dhabi", " September 23 2020", "Mumbai won by 49 runs", "Mumbai Indians", "Kolkata Knight Riders", "Krunal Pandya ", "1", "3", "0", "0", "33.33"),
    @(" Abu Dhabi", " October 11 2020", "Mumbai won by 5 wickets (with 2 balls remaining)", "Mumbai Indians", "Delhi Capitals", "Krunal Pandya ", "12", "7", "2", "0", "171.42"),
    @(" Sharjah", " October 04 2020", "Mumbai won by 34 runs", "Mumbai Indians", "Sunrisers Hyderabad", "Krunal Pandya ", "20", "4", "2", "2", "500.00")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = $r + 2
    $rowRange = $ws.Range("A" + $rowNum + ":K" + $rowNum)
    # Force Text storage so numeric-looking strings ("3", "100.00", ...) are not
    # auto-converted to numbers, then restore the default "Normal" style so no
    # lingering custom number format is left on the cells.
    $rowRange.NumberFormat = "@"
    for ($c = 0; $c -lt $data[$r].Length; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $data[$r][$c]
    }
    $rowRange.Style = "Normal"
}

